$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("A2").Value = "namen"
$ws.Range("A3").Value = "namekJ-uw>"
$ws.Range("A4").Value = "name?S0y sFT3"
$ws.Range("A5").Value = "name*"
$ws.Range("A6").Value = "nameuUPC"
